# Updated the dates in the excel workbook to match the RiverSMART Run Range.
# For every sheet, column A holds date values (serial numbers) starting at
# row 2 down to the last contiguous populated row. Shift each date forward
# by exactly one year (same month/day, year+1) to realign with the new
# RiverSMART run range.

$wb = $excel.ActiveWorkbook
$xlDown = -4121

$sheetNames = @("Reservoirs", "CoordinatedOps", "ForecastLocations", "InterveningFlow")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $lastRow = $ws.Cells.Item(1, 1).End($xlDown).Row

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $old = $cell.Value2
        if ($old -eq $null) {
            continue
        }
        $d = [datetime]::FromOADate($old)
        $d2 = $d.AddYears(1)
        $cell.Value2 = $d2.ToOADate()
    }
}

# --- View / selection state ---------------------------------------------

$wsRes = $wb.Worksheets.Item("Reservoirs")
$wsRes.Activate()
$wsRes.Range("A1:A1048576").Select()

$wsCo = $wb.Worksheets.Item("CoordinatedOps")
$wsCo.Activate()
$wsCo.Range("A1:A4").Select()

$wsFl = $wb.Worksheets.Item("ForecastLocations")
$wsFl.Activate()
$wsFl.Range("A1:A1048576").Select()

$wsIf = $wb.Worksheets.Item("InterveningFlow")
$wsIf.Activate()
$wsIf.Range("A1:A1048576").Select()

$wsRes.Activate()
